# fix alasan jika ditolak
# Adds two new vehicle ("Mobil") asset rows to the asset sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (new mobil asset #1) ---
# Set cells in the same order the source application inserts them
# (tanggal_inventarisir, no_rangka, no_mesin, tgl_bpkb, no_bpkb, no_polisi,
#  nama, merk, type) so new shared-string entries are created in the
# expected order, followed by the remaining columns.
$ws.Range("C4").Value = "26/11/2025"
$ws.Range("H4").Value = "MHG000SK021469"
$ws.Range("I4").Value = "GGE-1021673"
$ws.Range("M4").Value = "22/11/2008"
$ws.Range("N4").Value = "A3554895"
$ws.Range("L4").Value = "R 5572 HS"
$ws.Range("B4").Value = "Mobil"
$ws.Range("D4").Value = "Toyota"
$ws.Range("E4").Value = "Innova"
$ws.Range("A4").Value = "Tersedia"
$ws.Range("F4").Value = 2500
$ws.Range("G4").Value = "Hitam"
$ws.Range("J4").Value = 2020
$ws.Range("K4").Value = 2021
$ws.Range("O4").Value = 14000002
$ws.Range("P4").Value = "Bantuan PEMDA Banyumas"

# --- Row 5 (new mobil asset #2) ---
$ws.Range("C5").Value = "26/11/2026"
$ws.Range("H5").Value = "MHG000SK021470"
$ws.Range("I5").Value = "GGE-1021674"
$ws.Range("M5").Value = "22/11/2009"
$ws.Range("N5").Value = "A3554896"
$ws.Range("L5").Value = "R 5511 AB"
$ws.Range("B5").Value = "Mobil"
$ws.Range("D5").Value = "Honda"
$ws.Range("E5").Value = "Mobilio"
$ws.Range("A5").Value = "Tersedia"
$ws.Range("F5").Value = 2000
$ws.Range("G5").Value = "Hitam"
$ws.Range("J5").Value = 2020
$ws.Range("K5").Value = 2021
$ws.Range("O5").Value = 14000003
$ws.Range("P5").Value = "Bantuan PEMDA Banyumas"

# Move the active selection to match the author's final cursor position
$ws.Range("G5").Select()
